$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Correct the estoque_atualizado figures for the two existing sales (BEMOL S/A restock update) ---
$ws.Range("G2").Value = -39
$ws.Range("G3").Value = -110

# --- Append the new atypical sale row (id 4): 2025-06-11, BEMOL S/A, RING LIGHT ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-06-11"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = 2

$ws.Range("C4").Value = "BEMOL S/A"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "357392"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = 4408

$ws.Range("F4").Value = "RING LIGHT 10 POLEGADAS COM TRIPE"

$ws.Range("G4").Value = -444
$ws.Range("H4").Value = 1.02
$ws.Range("I4").Value = 0.16
